# Update burndown for sprint 2
#
# Sprint 2 actually-burned hours (column C, rows 8-12) were filled in for
# the first day of the sprint; C14 (Total Hours Left formula
# =SUM(C8:C12)) recalculates automatically as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 2 Burndown")

$ws.Range("C8").Value = 3
$ws.Range("C9").Value = 3
$ws.Range("C10").Value = 2
$ws.Range("C11").Value = 3
$ws.Range("C12").Value = 4

# Move the active selection, matching where the author left off editing.
$ws.Activate()
$ws.Range("D25").Select()

# Clean up the stale hidden "_xlchart" defined names that Excel had left
# behind from an earlier chart-source selection; they no longer point at
# anything meaningful.
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}
